# Update metadata date and the two "System URI" values in the Include sheets.

$wb = $excel.ActiveWorkbook

# Sheet "Metadata": row 8 holds the "Date" property -> update timestamp.
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# Sheet "Include #0": row 4 holds the "System URI" for TRE-R38-SpecialiteOrdinale.
$include0 = $wb.Worksheets.Item("Include #0")
$include0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R38-SpecialiteOrdinale"

# Sheet "Include #1": row 4 holds the "System URI" for TRE-R01-EnsembleSavoirFaire-CISIS.
$include1 = $wb.Worksheets.Item("Include #1")
$include1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R01-EnsembleSavoirFaire-CISIS"
